# Auto-generated Excel COM-interop script
# Syncs numeric leve-profit figures in Sheets/Maduin_Profits.xlsx per the scheduled runner's refreshed
# market-board pricing snapshot (currentAveragePrice* / LevePrice* / LeveProfit* columns, cols H:N).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 731.25
$ws.Range("I6").Value = 570.2
$ws.Range("J6").Value = 999.6667
$ws.Range("K6").Value = 1710.6
$ws.Range("L6").Value = 2999.0001
$ws.Range("M6").Value = -1598.6
$ws.Range("N6").Value = -3223.0001
# Row 8
$ws.Range("H8").Value = 105.85714
$ws.Range("I8").Value = 105.85714
$ws.Range("K8").Value = 317.57142
$ws.Range("M8").Value = -178.57142
# Row 21
$ws.Range("H21").Value = 500
$ws.Range("I21").Value = 500
$ws.Range("K21").Value = 500
$ws.Range("M21").Value = -32
# Row 23
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -266
# Row 31
$ws.Range("H31").Value = 92.8
$ws.Range("I31").Value = 78.5
$ws.Range("J31").Value = 150
$ws.Range("K31").Value = 235.5
$ws.Range("L31").Value = 450
$ws.Range("M31").Value = -5.5
$ws.Range("N31").Value = -910
# Row 39
$ws.Range("H39").Value = 2331.5
$ws.Range("I39").Value = 1499.5
$ws.Range("K39").Value = 4498.5
$ws.Range("M39").Value = -4202.5
# Row 58
$ws.Range("H58").Value = 233.33333
$ws.Range("J58").Value = 300
$ws.Range("L58").Value = 900
$ws.Range("N58").Value = -1200
# Row 64
$ws.Range("H64").Value = 12000
# Row 67
$ws.Range("H67").Value = 12000
# Row 103
$ws.Range("H103").Value = 3076.5557
$ws.Range("J103").Value = 3829.6667
$ws.Range("L103").Value = 11489.0001
$ws.Range("N103").Value = -12661.0001
# Row 116
$ws.Range("H116").Value = 2633.3333
$ws.Range("I116").Value = 2950
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2950
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 492
$ws.Range("N116").Value = -8884
# Row 132
$ws.Range("H132").Value = 4749.8335
$ws.Range("I132").Value = 5999.5
$ws.Range("K132").Value = 17998.5
$ws.Range("M132").Value = -15468.5
# Row 138
$ws.Range("H138").Value = 3952.0908
$ws.Range("I138").Value = 4137.25
$ws.Range("J138").Value = 3910.9443
$ws.Range("K138").Value = 12411.75
$ws.Range("L138").Value = 11732.8329
$ws.Range("M138").Value = -7271.75
$ws.Range("N138").Value = -22012.8329

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3887.3333
$ws.Range("I2").Value = 825
$ws.Range("K2").Value = 825
$ws.Range("M2").Value = -712
# Row 30
$ws.Range("H30").Value = 2488
$ws.Range("J30").Value = 2488
$ws.Range("L30").Value = 2488
$ws.Range("N30").Value = -2788
# Row 32
$ws.Range("H32").Value = 2947.577
$ws.Range("I32").Value = 2947.577
$ws.Range("K32").Value = 2947.577
$ws.Range("M32").Value = -2660.577
# Row 116
$ws.Range("H116").Value = 3887.3333
$ws.Range("I116").Value = 825
$ws.Range("K116").Value = 825
$ws.Range("M116").Value = 1469

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 408.625
$ws.Range("I22").Value = 383
$ws.Range("J22").Value = 424
$ws.Range("K22").Value = 383
$ws.Range("L22").Value = 424
$ws.Range("M22").Value = -210
$ws.Range("N22").Value = -770
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 94
$ws.Range("H94").Value = 2745.2
$ws.Range("I94").Value = 490.4
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 490.4
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -39.39999999999998
$ws.Range("N94").Value = -5902
# Row 107
$ws.Range("H107").Value = 576.75
$ws.Range("I107").Value = 512.8
$ws.Range("J107").Value = 683.3333
$ws.Range("K107").Value = 512.8
$ws.Range("L107").Value = 683.3333
$ws.Range("M107").Value = 1407.2
$ws.Range("N107").Value = -4523.3333
# Row 134
$ws.Range("H134").Value = 6416.3335
$ws.Range("I134").Value = 7166
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 21498
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -18963
$ws.Range("N134").Value = -22069.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 739
$ws.Range("I25").Value = 600
$ws.Range("J25").Value = 831.6667
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 831.6667
$ws.Range("M25").Value = -426
$ws.Range("N25").Value = -1179.6667
# Row 31
$ws.Range("H31").Value = 2867.7778
$ws.Range("I31").Value = 2101.375
$ws.Range("K31").Value = 2101.375
$ws.Range("M31").Value = -1806.375
# Row 34
$ws.Range("H34").Value = 2867.7778
$ws.Range("I34").Value = 2101.375
$ws.Range("K34").Value = 2101.375
$ws.Range("M34").Value = -1899.375
# Row 58
$ws.Range("H58").Value = 1699.25
$ws.Range("I58").Value = 1699.25
$ws.Range("K58").Value = 1699.25
$ws.Range("M58").Value = -1496.25
# Row 60
$ws.Range("H60").Value = 16379.4
$ws.Range("J60").Value = 27949
$ws.Range("L60").Value = 27949
$ws.Range("N60").Value = -28971
# Row 63
$ws.Range("H63").Value = 99994.5
$ws.Range("J63").Value = 99994.5
$ws.Range("L63").Value = 99994.5
$ws.Range("N63").Value = -101366.5
# Row 66
$ws.Range("H66").Value = 99994.5
$ws.Range("J66").Value = 99994.5
$ws.Range("L66").Value = 299983.5
$ws.Range("N66").Value = -306847.5
# Row 70
$ws.Range("H70").Value = 24999.5
$ws.Range("J70").Value = 24999.5
$ws.Range("L70").Value = 24999.5
$ws.Range("N70").Value = -25629.5
# Row 73
$ws.Range("H73").Value = 24999.5
$ws.Range("J73").Value = 24999.5
$ws.Range("L73").Value = 24999.5
$ws.Range("N73").Value = -27183.5
# Row 136
$ws.Range("H136").Value = 1699.25
$ws.Range("I136").Value = 1699.25
$ws.Range("K136").Value = 5097.75
$ws.Range("M136").Value = -2547.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 95
$ws.Range("H95").Value = 12999
$ws.Range("J95").Value = 12999
$ws.Range("L95").Value = 38997
$ws.Range("N95").Value = -43115
# Row 119
$ws.Range("H119").Value = 495
$ws.Range("I119").Value = 495
$ws.Range("K119").Value = 1485
$ws.Range("M119").Value = 3353
# Row 120
$ws.Range("H120").Value = 6250
$ws.Range("I120").Value = 2500
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 7500
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -2662
$ws.Range("N120").Value = -39676

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1750
$ws.Range("J80").Value = 1500
$ws.Range("L80").Value = 1500
$ws.Range("N80").Value = -3496
# Row 83
$ws.Range("H83").Value = 1750
$ws.Range("J83").Value = 1500
$ws.Range("L83").Value = 7500
$ws.Range("N83").Value = -17484
# Row 101
$ws.Range("H101").Value = 34352
$ws.Range("J101").Value = 34352
$ws.Range("L101").Value = 34352
$ws.Range("N101").Value = -40842
# Row 126
$ws.Range("H126").Value = 7568.4287
$ws.Range("J126").Value = 6666.6665
$ws.Range("L126").Value = 19999.9995
$ws.Range("N126").Value = -24939.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1905.8948
$ws.Range("I7").Value = 1694.8823
$ws.Range("K7").Value = 1694.8823
$ws.Range("M7").Value = -1582.8823
# Row 12
$ws.Range("H12").Value = 1445
$ws.Range("J12").Value = 1990
$ws.Range("L12").Value = 1990
$ws.Range("N12").Value = -2330
# Row 16
$ws.Range("H16").Value = 361.66666
$ws.Range("I16").Value = 361.66666
$ws.Range("K16").Value = 361.66666
$ws.Range("M16").Value = -191.66666
# Row 40
$ws.Range("H40").Value = 1413.5555
$ws.Range("I40").Value = 1413.5555
$ws.Range("K40").Value = 1413.5555
$ws.Range("M40").Value = -1277.5555
# Row 46
$ws.Range("H46").Value = 3112.0952
$ws.Range("I46").Value = 2362.5
$ws.Range("J46").Value = 4111.5557
$ws.Range("K46").Value = 2362.5
$ws.Range("L46").Value = 4111.5557
$ws.Range("M46").Value = -2174.5
$ws.Range("N46").Value = -4487.5557
# Row 122
$ws.Range("H122").Value = 5004.0586
$ws.Range("I122").Value = 4804.6
$ws.Range("K122").Value = 14413.8
$ws.Range("M122").Value = -11963.8
# Row 126
$ws.Range("H126").Value = 1905.8948
$ws.Range("I126").Value = 1694.8823
$ws.Range("K126").Value = 5084.6469
$ws.Range("M126").Value = -2614.6469

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 57
$ws.Range("H57").Value = 45000
$ws.Range("J57").Value = 45000
$ws.Range("L57").Value = 45000
$ws.Range("N57").Value = -46508
# Row 87
$ws.Range("H87").Value = 47555
$ws.Range("J87").Value = 47555
$ws.Range("L87").Value = 47555
$ws.Range("N87").Value = -50051
# Row 90
$ws.Range("H90").Value = 47555
$ws.Range("J90").Value = 47555
$ws.Range("L90").Value = 142665
$ws.Range("N90").Value = -155145
# Row 100
$ws.Range("H100").Value = 5363174.5
$ws.Range("I100").Value = 17424968
$ws.Range("J100").Value = 2377.6667
$ws.Range("K100").Value = 34849936
$ws.Range("L100").Value = 4755.3334
$ws.Range("M100").Value = -34849395
$ws.Range("N100").Value = -5837.3334
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 919.8461
$ws.Range("I126").Value = 958.5454999999999
$ws.Range("J126").Value = 707
$ws.Range("K126").Value = 2875.6365
$ws.Range("L126").Value = 2121
$ws.Range("M126").Value = -405.6364999999996
$ws.Range("N126").Value = -7061

